# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates the "K" column (column G) values on Sheet1 for rows 2-71 to
# reflect the newly-calculated strike counts (previously using a different
# statistic, now derived from "K" per the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 0
    6  = 2
    7  = 3
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 1
    14 = 2
    15 = 0
    16 = 2
    17 = 2
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 2
    24 = 1
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 3
    32 = 0
    33 = 1
    35 = 2
    37 = 2
    38 = 3
    39 = 1
    40 = 3
    41 = 0
    42 = 1
    43 = 1
    44 = 2
    45 = 1
    46 = 1
    47 = 2
    48 = 1
    49 = 2
    50 = 2
    51 = 0
    52 = 3
    53 = 0
    54 = 2
    55 = 1
    56 = 2
    57 = 2
    58 = 1
    59 = 1
    60 = 1
    61 = 1
    62 = 2
    63 = 2
    64 = 1
    65 = 1
    66 = 2
    67 = 1
    68 = 1
    70 = 2
    71 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
